$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap values between row 2 and row 3 for columns D, J, K, L, M, P
$ws.Range("D2").Value = 44875
$ws.Range("J2").Value = 1000
$ws.Range("K2").Value = 1600
$ws.Range("L2").Value = 1700
$ws.Range("M2").Value = 1650
$ws.Range("P2").Value = 1650

$ws.Range("D3").Value = 44547
$ws.Range("J3").Value = 400
$ws.Range("K3").Value = 1500
$ws.Range("L3").Value = 1600
$ws.Range("M3").Value = 1550
$ws.Range("P3").Value = 1550
